# Apply changes described by the commit:
# "changed MP time limit and corrected error in fixed recourse data"
#
# 1. The Status value "OPTIMAL" (shared across all data rows) becomes "TIME_LIMIT".
# 2. The objective (B), gap (C) and solve time (D) values for rows 2-11 are corrected.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# 1) Update Status column (E2:E11) from OPTIMAL to TIME_LIMIT
$ws.Range("E2:E11").Value = "TIME_LIMIT"

# 2) Update objective / gap / solve time values for each instance row
$data = @{
    2  = @{ B = -411.51897100078946; C = 7.896805345891741;  D = 3604.548796775 }
    3  = @{ B = -411.3960477688105;  C = 6.492730181926509;  D = 3652.929948378 }
    4  = @{ B = -412.3476210996386;  C = 5.181024663335519;  D = 3799.723364972 }
    5  = @{ B = -416.71129461895856; C = 2.627468796449248;  D = 3858.943483408 }
    6  = @{ B = -408.2905883517202;  C = 6.704542726035493;  D = 3600.528406784 }
    7  = @{ B = -403.99029359095584; C = 4.17639584236483;   D = 3768.484232653 }
    8  = @{ B = -400.19199641565194; C = 1.1297632678732288; D = 3832.915016989 }
    9  = @{ B = -412.09265615224575; C = 6.122608144824806;  D = 3663.325041742 }
    10 = @{ B = -408.4459810048513;  C = 6.617595411152933;  D = 3627.94769365  }
    11 = @{ B = -403.1039288610365;  C = 1.069073315486184;  D = 3812.042297833 }
}

foreach ($row in $data.Keys) {
    $vals = $data[$row]
    $ws.Range("B$row").Value = $vals.B
    $ws.Range("C$row").Value = $vals.C
    $ws.Range("D$row").Value = $vals.D
}
